$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Adatok")

# Dates in column G are stored as plain text (e.g. "2023-05-10"), not real
# date serials. Force the new date cells to Text format first so Excel
# doesn't auto-convert the "yyyy-mm-dd" looking strings into date values.
$ws.Range("G7:G10").NumberFormat = "@"

# Fix row 7: replace the old "probakép" placeholder entry with the real
# "Kantaros farmer" loan record.
$ws.Range("A7").Value = "Ruha"
$ws.Range("B7").Value = "Kantáros farmer "
$ws.Range("C7").Value = 86
$ws.Range("D7").Value = "Kantáros farmer ruha csíkos bodyval, új"
$ws.Range("E7").Value = "nem"
$ws.Range("F7").Value = "Saját"
$ws.Range("G7").Value = "2023-05-15"

# Add new rows 8-10 for the additional loaned items.
$ws.Range("A8").Value = "Felső"
$ws.Range("B8").Value = "Rózsaszín szett"
$ws.Range("C8").Value = 80
$ws.Range("D8").Value = "Rózsaszín átlapolt felső nadrággal és sapkával, új"
$ws.Range("E8").Value = "nem"
$ws.Range("F8").Value = "Saját"
$ws.Range("G8").Value = "2023-05-15"

$ws.Range("A9").Value = "Egyéb"
$ws.Range("B9").Value = "Mályva overál"
$ws.Range("C9").Value = 68
$ws.Range("D9").Value = "Mályva színű,bundás overál"
$ws.Range("E9").Value = "igen"
$ws.Range("F9").Value = "Dóri, Krisztián"
$ws.Range("G9").Value = "2023-05-15"

$ws.Range("A10").Value = "Hálózsák"
$ws.Range("B10").Value = "Vonatos hálózsák"
$ws.Range("C10").Value = 65
$ws.Range("D10").Value = "Vastag, vonat mintás hálózsák"
$ws.Range("E10").Value = "igen"
$ws.Range("F10").Value = "Dóri, Krisztián"
$ws.Range("G10").Value = "2023-05-15"
